# Apply 2025-02-10 data update to violent-crime-full-year workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7900
$ws.Range("L2").Value = 610
$ws.Range("L3").Value = 612
$ws.Range("J4").Value = 1852
$ws.Range("K4").Value = 1735
$ws.Range("L4").Value = 161
$ws.Range("L5").Value = 49
$ws.Range("L6").Value = 681
$ws.Range("J7").Value = 29321
$ws.Range("K7").Value = 27526
$ws.Range("L7").Value = 2113

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 33
$ws.Range("L3").Value = 36
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 120

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 22
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 73

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 75
$ws.Range("L8").Value = 120
$ws.Range("L18").Value = 15
$ws.Range("L19").Value = 71
$ws.Range("L21").Value = 6
$ws.Range("K24").Value = 88
$ws.Range("L24").Value = 2
$ws.Range("L29").Value = 108
$ws.Range("L33").Value = 87
$ws.Range("L37").Value = 73
$ws.Range("L42").Value = 73
$ws.Range("L44").Value = 14
$ws.Range("L48").Value = 35
$ws.Range("L51").Value = 30
$ws.Range("L52").Value = 42
$ws.Range("L53").Value = 27
$ws.Range("L55").Value = 22
$ws.Range("L60").Value = 15
$ws.Range("J63").Value = 200
$ws.Range("K63").Value = 76
$ws.Range("L63").Value = 8
$ws.Range("L65").Value = 41
$ws.Range("L67").Value = 67
$ws.Range("L75").Value = 9
$ws.Range("L76").Value = 29
$ws.Range("L78").Value = 21
$ws.Range("L79").Value = 61
$ws.Range("L85").Value = 103
$ws.Range("L88").Value = 34
$ws.Range("L89").Value = 23
$ws.Range("L93").Value = 12
$ws.Range("J101").Value = 29321
$ws.Range("K101").Value = 27526
$ws.Range("L101").Value = 2113

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 19
$ws.Range("L4").Value = 6
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L3").Value = 33
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 73

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K4").Value = 9
$ws.Range("L4").Value = 1
$ws.Range("K7").Value = 88
$ws.Range("L7").Value = 2

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("L3").Value = 1
$ws.Range("L7").Value = 6

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 27
$ws.Range("L5").Value = 3
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 15
$ws.Range("L3").Value = 11
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 42
